# Generate Report for Archive
#
# Updates the localization status of the two files that have moved from
# "Ready for handoff" into "In Translation":
#   - 31e3b9b5-c1f0-4455-8e71-e8931a544c34.md
#   - 7c05bd8b-9d38-4130-971e-e4d4e802e8de.md
#
# This status is reflected in three places:
#   - the "Overview" sheet (zh-cn and de-de columns for the matching rows)
#   - the "zh-cn" sheet (Status column for the matching rows)
#   - the "de-de" sheet (Status column for the matching rows)

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
